$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark that currently sits after the paragraph
#    "Met eigen laptop op gastnetwerk komen".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Recolor the "Actiebeschrijving meegeven aan volgende kamer (...)"
#    bullet from orange (FFC000) to green (00B050).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Actiebeschrijving meegeven aan volgende kamer*") {
        $p.Range.Font.Color = 5287936   # RGB(0x00,0xB0,0x50) == 00B050
    }
}

# ---------------------------------------------------------------------------
# 3) Give the "Kamerteller" bullet the same green (00B050) color.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Kamerteller*") {
        $p.Range.Font.Color = 5287936   # RGB(0x00,0xB0,0x50) == 00B050
    }
}

# ---------------------------------------------------------------------------
# 4) Drop the "Meelopen met iemand door een deur" bullet entirely.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Meelopen met iemand door een deur*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 5) Re-add the "_GoBack" bookmark, now at the very start of the
#    "Spelbegin: keuze tussen direct beginnen, of uitleg" paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Spelbegin: keuze tussen direct beginnen*") {
        $ins = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $ins)
        break
    }
}

# ---------------------------------------------------------------------------
# 6) Move the "lastRenderedPageBreak" marker: it currently sits in the run
#    "Statusdeel op scherm" and should instead sit in the run "I" (the
#    first letter of "Inventory").
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Statusdeel op scherm*") {
        $pStart = $p.Range.Start
        $runRange = $d.Range($pStart, $pStart + 20)   # "Statusdeel op scherm"
        $runRange.Delete()
        $insPoint = $d.Range($pStart, $pStart)
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
               '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body><w:p><w:r w:rsidRPr="00EF02F3"><w:rPr><w:color w:val="00B050"/></w:rPr>' +
               '<w:t>Statusdeel op scherm</w:t></w:r></w:p></w:body></w:document>' +
               '</pkg:xmlData></pkg:part></pkg:package>'
        $insPoint.InsertXML($xml)
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Inventory*") {
        $pStart = $p.Range.Start
        $insPoint = $d.Range($pStart, $pStart)
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
               '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body><w:p><w:r w:rsidRPr="00EF02F3"><w:rPr><w:color w:val="00B050"/></w:rPr>' +
               '<w:lastRenderedPageBreak/><w:t>I</w:t></w:r></w:p></w:body></w:document>' +
               '</pkg:xmlData></pkg:part></pkg:package>'
        $insPoint.InsertXML($xml)
        $oldI = $d.Range($pStart + 1, $pStart + 2)
        $oldI.Delete()
        break
    }
}
